# mBom.xlsx update: add 4 new parts (rows 33-36) to the BOM table,
# adjust row 32's height, and move the view/selection to the new bottom
# of the sheet (matches the upstream commit "updated bom and apparently
# gantt chart").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("mBom")

# ---------------------------------------------------------------------
# Row 32 already exists - only its height changes (14.9 -> 14.3)
# ---------------------------------------------------------------------
$ws.Rows.Item(32).RowHeight = 14.3

# ---------------------------------------------------------------------
# Helper: clone the formatting of an existing row-32 cell onto a new
# cell before writing its value, so the new rows pick up the same cell
# styles already used in the sheet (row 32 alone carries every style we
# need: style 14 from B32/G32/H32, style 15 from F32/I32, style 1 from
# D32/E32, style 2 from C32, style 3 from K32:N32).
# ---------------------------------------------------------------------
function Copy-Format($srcAddr, $dstAddr) {
    $ws.Range($srcAddr).Copy()
    $ws.Range($dstAddr).PasteSpecial(-4122)
}

# =======================================================================
# Row 33 - 100 Ohm Resistor
# =======================================================================
Copy-Format "F32" "B33"
Copy-Format "C32" "C33"
Copy-Format "D32" "D33"
Copy-Format "D32" "E33"
Copy-Format "F32" "F33"
Copy-Format "D32" "G33"
Copy-Format "B32" "H33"
Copy-Format "D32" "I33"
Copy-Format "K32" "K33"
Copy-Format "K32" "L33"
Copy-Format "K32" "M33"
Copy-Format "K32" "N33"

$ws.Rows.Item(33).RowHeight = 14.9
$ws.Cells.Item(33, 2).Value = "100 Ohm Resistor"
$ws.Cells.Item(33, 3).Value = "0603"
$ws.Cells.Item(33, 4).Value = 100
$ws.Cells.Item(33, 5).Value = "meh "
$ws.Cells.Item(33, 6).Value = "RES SMD 100 OHM 1% 1/8W 0603 "
$ws.Cells.Item(33, 7).Value = "Vishay Beyschlag"
$ws.Cells.Item(33, 8).Value = "MCT06030C1000FP500"
$ws.Cells.Item(33, 9).Value = "MCT0603-100-CFCT-ND"
$ws.Cells.Item(33, 11).Value = 0.08
$ws.Cells.Item(33, 12).Value = 0.069
$ws.Cells.Item(33, 13).Value = 0.0408
$ws.Cells.Item(33, 14).Value = 0.024

# =======================================================================
# Row 34 - 649 Ohm Resistor
# =======================================================================
Copy-Format "D32" "B34"
Copy-Format "C32" "C34"
Copy-Format "D32" "D34"
Copy-Format "F32" "F34"
Copy-Format "D32" "G34"
Copy-Format "B32" "H34"
Copy-Format "F32" "I34"
Copy-Format "K32" "K34"
Copy-Format "K32" "L34"
Copy-Format "K32" "M34"
Copy-Format "K32" "N34"

$ws.Rows.Item(34).RowHeight = 28.35
$ws.Cells.Item(34, 2).Value = "649 Ohm Resistor"
$ws.Cells.Item(34, 3).Value = "0603"
$ws.Cells.Item(34, 4).Value = 649
$ws.Cells.Item(34, 6).Value = "RES SMD 649 OHM 1% 1/10W 0603 "
$ws.Cells.Item(34, 7).Value = "Stackpole Electronics Incx"
$ws.Cells.Item(34, 8).Value = "RMCF0603FT649R"
$ws.Cells.Item(34, 9).Value = "RMCF0603FT649RCT-ND "
$ws.Cells.Item(34, 11).Value = 0.1
$ws.Cells.Item(34, 12).Value = 0.025
$ws.Cells.Item(34, 13).Value = 0.0105
$ws.Cells.Item(34, 14).Value = 0.00384

# =======================================================================
# Row 35 - 59k Ohm Resistor
# =======================================================================
Copy-Format "D32" "B35"
Copy-Format "C32" "C35"
Copy-Format "D32" "D35"
Copy-Format "F32" "F35"
Copy-Format "D32" "G35"
Copy-Format "B32" "H35"
Copy-Format "F32" "I35"
Copy-Format "K32" "K35"
Copy-Format "K32" "L35"
Copy-Format "K32" "M35"
Copy-Format "K32" "N35"

$ws.Rows.Item(35).RowHeight = 14.9
$ws.Cells.Item(35, 2).Value = "59k Ohm Resistor"
$ws.Cells.Item(35, 3).Value = "0603"
$ws.Cells.Item(35, 4).Value = "59k"
$ws.Cells.Item(35, 6).Value = "RES SMD 59K OHM 1% 1/10W 0603 "
$ws.Cells.Item(35, 7).Value = "Panasonic electronic Components"
$ws.Cells.Item(35, 8).Value = "ERJ-3EKF5902V"
$ws.Cells.Item(35, 9).Value = "P59.0KHCT-ND "
$ws.Cells.Item(35, 11).Value = 0.1
$ws.Cells.Item(35, 12).Value = 0.015
$ws.Cells.Item(35, 13).Value = 0.0114
$ws.Cells.Item(35, 14).Value = 0.00416

# =======================================================================
# Row 36 - 10pF Capacitor
# =======================================================================
Copy-Format "D32" "B36"
Copy-Format "C32" "C36"
Copy-Format "D32" "D36"
Copy-Format "D32" "E36"
Copy-Format "F32" "F36"
Copy-Format "D32" "G36"
Copy-Format "D32" "H36"
Copy-Format "F32" "I36"
Copy-Format "K32" "K36"
Copy-Format "K32" "L36"
Copy-Format "K32" "M36"
Copy-Format "K32" "N36"

$ws.Rows.Item(36).RowHeight = 14.9
$ws.Cells.Item(36, 2).Value = "10pF Capacitor"
$ws.Cells.Item(36, 3).Value = "0402"
$ws.Cells.Item(36, 4).Value = "10pf"
$ws.Cells.Item(36, 5).Value = "50V"
$ws.Cells.Item(36, 6).Value = "CAP CER 10PF 50V C0G 0402 "
$ws.Cells.Item(36, 7).Value = "TDK Corporation"
$ws.Cells.Item(36, 8).Value = "C1005C0G1H100C050BA"
$ws.Cells.Item(36, 9).Value = "445-4896-1-ND "
$ws.Cells.Item(36, 11).Value = 0.1
$ws.Cells.Item(36, 12).Value = 0.026
$ws.Cells.Item(36, 13).Value = 0.011
$ws.Cells.Item(36, 14).Value = 0.0063

# ---------------------------------------------------------------------
# Scroll the view down to the new bottom of the table and match the
# recorded selection/active cell.
# ---------------------------------------------------------------------
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
$ws.Range("N37").Select() | Out-Null
